# Updated cryptos list on Sat Jul 13 05:35:25 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "57.954.44"
$ws.Range("E2").Value = "  +1.52%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.125.96"
$ws.Range("E3").Value = "  +1.29%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 - BNB (numeric-looking -> force text with leading apostrophe)
$ws.Range("D5").Value = "'534.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.59%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'138.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.27%  "

# Row 8 - XRP
$ws.Range("D8").Value = "'0.508"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +12.21%  "

# Row 9 - Toncoin
$ws.Range("D9").Value = "'7.35"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.07%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "'0.109"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.13%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  +4.67%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +3.28%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.662.19"
$ws.Range("E13").Value = "  +1.26%  "

# Row 14 - Avalanche
$ws.Range("E14").Value = "  +1.41%  "

# Row 15 - ShibaInu
$ws.Range("D15").Value = "'0.0000169"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.49%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "58.030.60"
$ws.Range("E16").Value = "  +1.53%  "

# Row 17 - now Polkadot (was WrappedEther)
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").Value = "'6.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.38%  "

# Row 18 - now WrappedEther (was Polkadot)
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.127.08"
$ws.Range("E18").Value = "  +1.33%  "

# Row 19 - Chainlink
$ws.Range("E19").Value = "  +3.06%  "

# Row 20 - Uniswap
$ws.Range("D20").Value = "'8.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.66%  "

# Row 21 - BitcoinCash
$ws.Range("D21").Value = "'376.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +8.17%  "

# Row 22 - Dai
$ws.Range("D22").Value = "'0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.06%  "

# Row 23 - LEO
$ws.Range("D23").Value = "'5.73"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.74%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "'69.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.27%  "

# Row 25 - Polygon
$ws.Range("D25").Value = "'0.511"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.44%  "

# Row 26 - Kaspa
$ws.Range("E26").Value = "  +0.67%  "

# Row 27 - Binance-PegBSC-USD
$ws.Range("E27").Value = "  -0.01%  "

# Row 28 - PEPE (contains unicode subscript three, not numeric-parseable anyway)
$sub3 = [char]0x2083
$ws.Range("D28").Value = "{0}{1}{2}" -f "0.0", $sub3, "0888"
$ws.Range("E28").Value = "  +2.61%  "

# Row 29 - InternetComputer(DFINITY)
$ws.Range("D29").Value = "'7.73"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.80%  "

# Row 30 - RenderToken
$ws.Range("D30").Value = "'6.16"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.21%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +0.58%  "

# Row 32 - EthereumClassic
$ws.Range("D32").Value = "'21.69"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.08%  "

# Row 33 - NEARProtocol
$ws.Range("E33").Value = "  +4.89%  "

# Row 34 - Fetch.AI
$ws.Range("E34").Value = "  +3.06%  "

# Row 35 - Monero
$ws.Range("D35").Value = "'160.99"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.90%  "

# Row 36 - Aptos
$ws.Range("D36").Value = "'6.21"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.06%  "

# Row 37 - ImmutableX
$ws.Range("E37").Value = "  +6.89%  "

# Row 38 - EnergySwap
$ws.Range("D38").Value = "'25.55"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.32%  "

# Row 39 - Stacks
$ws.Range("D39").Value = "'1.66"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.51%  "

# Row 40 - Hedera
$ws.Range("D40").Value = "'0.0675"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.94%  "

# Row 41 - Filecoin
$ws.Range("E41").Value = "  +3.76%  "

# Row 42 - Maker
$ws.Range("D42").Value = "2.551.21"
$ws.Range("E42").Value = "  +6.68%  "

# Row 43 - OKB
$ws.Range("D43").Value = "'38.70"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.63%  "

# Row 44 - Mantle
$ws.Range("D44").Value = "'0.699"
$ws.Range("D44").Style = "Normal"

# Row 45 - VeChain
$ws.Range("D45").Value = "'0.0272"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.86%  "

# Row 46 - FirstDigitalUSD
$ws.Range("E46").Value = "  +0.00%  "

# Row 47 - Cosmos
$ws.Range("D47").Value = "'6.21"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.31%  "

# Row 48 - ONDO
$ws.Range("D48").Value = "'0.978"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.88%  "

# Row 49 - Stellar
$ws.Range("D49").Value = "'0.0984"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.43%  "

# Row 50 - InjectiveProtocol
$ws.Range("D50").Value = "'20.02"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.64%  "

# Row 51 - SuiNetwork
$ws.Range("D51").Value = "'0.749"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.66%  "
